$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff.
# D/E columns get NumberFormat '@' (Text) applied first so that
# numeric-looking strings (e.g. '573.13', '0.0000170') are preserved
# verbatim as text instead of being coerced into floating point numbers
# by Excel's automatic type inference.
$updates = @(
    @{ Cell = 'D2'; Value = '61.082.92' }
    @{ Cell = 'E2'; Value = '  -1.92%  ' }
    @{ Cell = 'D3'; Value = '2.433.89' }
    @{ Cell = 'E3'; Value = '  -0.45%  ' }
    @{ Cell = 'E4'; Value = '  -0.13%  ' }
    @{ Cell = 'D5'; Value = '573.13' }
    @{ Cell = 'E5'; Value = '  -1.48%  ' }
    @{ Cell = 'D6'; Value = '140.59' }
    @{ Cell = 'E6'; Value = '  -2.12%  ' }
    @{ Cell = 'E7'; Value = '  +0.16%  ' }
    @{ Cell = 'D9'; Value = '2.420.90' }
    @{ Cell = 'E9'; Value = '  -0.91%  ' }
    @{ Cell = 'E10'; Value = '  +1.31%  ' }
    @{ Cell = 'E11'; Value = '  +0.88%  ' }
    @{ Cell = 'E12'; Value = '  -1.49%  ' }
    @{ Cell = 'E13'; Value = '  -1.51%  ' }
    @{ Cell = 'D14'; Value = '26.11' }
    @{ Cell = 'E14'; Value = '  -1.23%  ' }
    @{ Cell = 'D15'; Value = '0.0000170' }
    @{ Cell = 'D16'; Value = '2.860.41' }
    @{ Cell = 'D17'; Value = '61.013.27' }
    @{ Cell = 'E17'; Value = '  -1.79%  ' }
    @{ Cell = 'D18'; Value = '2.422.66' }
    @{ Cell = 'E18'; Value = '  -0.68%  ' }
    @{ Cell = 'D19'; Value = '10.59' }
    @{ Cell = 'E19'; Value = '  -2.94%  ' }
    @{ Cell = 'D20'; Value = '7.30' }
    @{ Cell = 'E20'; Value = '  +2.54%  ' }
    @{ Cell = 'D21'; Value = '324.03' }
    @{ Cell = 'E21'; Value = '  -2.01%  ' }
    @{ Cell = 'D22'; Value = '4.05' }
    @{ Cell = 'E22'; Value = '  -1.52%  ' }
    @{ Cell = 'D23'; Value = '6.15' }
    @{ Cell = 'E23'; Value = '  +2.61%  ' }
    @{ Cell = 'E24'; Value = '  +0.08%  ' }
    @{ Cell = 'E25'; Value = '  -3.83%  ' }
    @{ Cell = 'D26'; Value = '65.23' }
    @{ Cell = 'E26'; Value = '  -0.95%  ' }
    @{ Cell = 'D27'; Value = '8.88' }
    @{ Cell = 'E27'; Value = '  -5.28%  ' }
    @{ Cell = 'D28'; Value = '580.37' }
    @{ Cell = 'E28'; Value = '  -6.27%  ' }
    @{ Cell = 'D29'; Value = '2.558.68' }
    @{ Cell = 'E29'; Value = '  -0.23%  ' }
    @{ Cell = 'B30'; Value = 'PEPE' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe' }
    @{ Cell = 'D30'; Value = '0.0₃0916' }
    @{ Cell = 'E30'; Value = '  -4.11%  ' }
    @{ Cell = 'B31'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D31'; Value = '7.90' }
    @{ Cell = 'E31'; Value = '  -1.34%  ' }
    @{ Cell = 'B32'; Value = 'Fetch.AI' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' }
    @{ Cell = 'D32'; Value = '1.35' }
    @{ Cell = 'E32'; Value = '  -5.36%  ' }
    @{ Cell = 'B33'; Value = 'PancakeSwap' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Cell = 'D33'; Value = '1.84' }
    @{ Cell = 'E33'; Value = '  -1.97%  ' }
    @{ Cell = 'B34'; Value = 'Kaspa' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ Cell = 'D34'; Value = '0.133' }
    @{ Cell = 'E34'; Value = '  -6.07%  ' }
    @{ Cell = 'B35'; Value = 'FirstDigitalUSD' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd' }
    @{ Cell = 'D35'; Value = '1.00' }
    @{ Cell = 'E35'; Value = '  +0.15%  ' }
    @{ Cell = 'B36'; Value = 'NEARProtocol' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Cell = 'D36'; Value = '4.62' }
    @{ Cell = 'E36'; Value = '  -5.97%  ' }
    @{ Cell = 'B37'; Value = 'PolygonEcosystemToken' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol' }
    @{ Cell = 'D37'; Value = '0.369' }
    @{ Cell = 'E37'; Value = '  -1.73%  ' }
    @{ Cell = 'B38'; Value = 'Monero' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D38'; Value = '150.69' }
    @{ Cell = 'E38'; Value = '  -0.43%  ' }
    @{ Cell = 'B39'; Value = 'ImmutableX' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D39'; Value = '1.38' }
    @{ Cell = 'E39'; Value = '  -3.85%  ' }
    @{ Cell = 'B40'; Value = 'EthereumClassic' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Cell = 'D40'; Value = '18.24' }
    @{ Cell = 'E40'; Value = '  -0.53%  ' }
    @{ Cell = 'B41'; Value = 'RenderToken' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render' }
    @{ Cell = 'D41'; Value = '5.12' }
    @{ Cell = 'E41'; Value = '  -2.37%  ' }
    @{ Cell = 'B42'; Value = 'USDe' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde' }
    @{ Cell = 'D42'; Value = '0.999' }
    @{ Cell = 'E42'; Value = '  +0.05%  ' }
    @{ Cell = 'B43'; Value = 'OKB' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D43'; Value = '41.72' }
    @{ Cell = 'E43'; Value = '  -1.98%  ' }
    @{ Cell = 'B44'; Value = 'Stacks' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Cell = 'D44'; Value = '1.66' }
    @{ Cell = 'E44'; Value = '  -5.92%  ' }
    @{ Cell = 'B45'; Value = 'dogwifhat' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' }
    @{ Cell = 'D45'; Value = '2.35' }
    @{ Cell = 'E45'; Value = '  -4.62%  ' }
    @{ Cell = 'B46'; Value = 'BabyDogeCoin' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge' }
    @{ Cell = 'D46'; Value = '0.0₆0286' }
    @{ Cell = 'E46'; Value = '  +25.97%  ' }
    @{ Cell = 'B47'; Value = 'Aave' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = 'D47'; Value = '141.65' }
    @{ Cell = 'E47'; Value = '  -1.25%  ' }
    @{ Cell = 'B48'; Value = 'Filecoin' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = 'D48'; Value = '3.53' }
    @{ Cell = 'E48'; Value = '  -2.71%  ' }
    @{ Cell = 'B49'; Value = 'Mantle' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' }
    @{ Cell = 'D49'; Value = '0.594' }
    @{ Cell = 'E49'; Value = '  -0.74%  ' }
    @{ Cell = 'B50'; Value = 'InjectiveProtocol' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Cell = 'D50'; Value = '19.62' }
    @{ Cell = 'E50'; Value = '  +0.52%  ' }
    @{ Cell = 'B51'; Value = 'Hedera' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D51'; Value = '0.0507' }
    @{ Cell = 'E51'; Value = '  -3.51%  ' }
)

foreach ($u in $updates) {
    $col = $u.Cell.Substring(0,1)
    $range = $ws.Range($u.Cell)
    if ($col -eq 'D' -or $col -eq 'E') {
        $range.NumberFormat = '@'
    }
    $range.Value = $u.Value
}
